$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 401, shifting existing rows 401-459 down to 402-460
$ws.Rows("401").Insert()

# Populate the newly inserted row 401 with the new record
$ws.Range("A401").Value = 4
$ws.Range("B401").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C401").Value = "Los Lagos"
$ws.Range("D401").Value2 = 45077
$ws.Range("D401").NumberFormat = $ws.Range("D400").NumberFormat
$ws.Range("E401").Value = 10
$ws.Range("F401").Value = 100112037
$ws.Range("G401").Value = "Cebollín"
$ws.Range("H401").Value = "Sin especificar"
$ws.Range("I401").Value = "Primera"
$ws.Range("J401").Value = 35
$ws.Range("K401").Value = 6500
$ws.Range("L401").Value = 6500
$ws.Range("M401").Value = 6500
$ws.Range("N401").Value = "$/paquete 36 unidades"
$ws.Range("O401").Value = "Región Metropolitana"
$ws.Range("P401").Value = 181
$ws.Range("Q401").Value = 36
$ws.Range("R401").Value = "Hortaliza"
